# Add a new "mm" data column (F) to the dataSources sheet, mirroring the
# existing es/nn/ads/city columns, then leave the selection on F10 as the
# author did after finishing the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataSources")

# Bring in the new values, matching the formatting already used by the
# neighboring "city" column (E) so the new column looks consistent.
$ws.Range("E1:E6").Copy($ws.Range("F1:F6"))

$ws.Range("F1").Value = "mm1"
$ws.Range("F2").Value = "mm2"
$ws.Range("F3").Value = "mm3"
$ws.Range("F4").Value = "mm4"
$ws.Range("F5").Value = "mm5"
$ws.Range("F6").Value = "mm6"

$ws.Range("F10").Select()
